# IC list for ctrl.
#
# The workbook originally had two tabs: an empty "Sheet2" (first tab) and a
# "Sheet1" (second tab, the active one) holding the actual 8-bit-CPU IC /
# wiring data. This edit removes the unused empty sheet and extends the
# "V#" placeholder rows at the bottom of the remaining sheet into a real
# per-net IC lookup table (reg / ram / pc / alu / ctrl columns), which is
# why the previously-unused "V27".."V63" placeholder strings disappear.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the blank "Sheet2" tab - all real content lives on "Sheet1".
$blank = $wb.Worksheets.Item("Sheet2")
[void]$blank.Delete()

$ws = $wb.Worksheets.Item("Sheet1")

# Wipe out the old placeholder rows (25 was a stray row, 26-87 held the
# "V1".."V63" placeholder list) so no stale values leak into the rebuilt
# table below.
$ws.Range("A25:S87").Clear()

# New section header (row 26): category columns for the per-net IC table.
$ws.Range("A26").Value = "CTRL"
$ws.Range("B26").Value = "reg"
$ws.Range("C26").Value = "ram"
$ws.Range("D26").Value = "pc"
$ws.Range("E26").Value = "alu"
$ws.Range("F26").Value = "ctrl"

# Per-net ("V1".."V26") IC part numbers used for each category.
# Columns: label, reg, ram, pc, alu, ctrl
$data = @()
$data += ,@("V1", 7411, 7430, 7430, 7408, 7400)
$data += ,@("V2", 7410, 7411, 7430, 7400, 7410)
$data += ,@("V3", 7410, 7410, 7430, 7400, 7430)
$data += ,@("V4", 7404, 7432, 7400, 7410, 74245)
$data += ,@("V5", 7410, 7404, 7432, 7400, 74245)
$data += ,@("V6", 7408, 7430, 7408, 7400, 74245)
$data += ,@("V7", 7430, 7410, 7432, 7411, $null)
$data += ,@("V8", 7430, 7400, 7400, 7400, $null)
$data += ,@("V9", 7400, 7400, 7408, 7410, $null)
$data += ,@("V10", 7432, 7410, 7410, 7400, $null)
$data += ,@("V11", 7400, 7432, 7432, 7410, $null)
$data += ,@("V12", $null, $null, 7404, 7410, $null)
$data += ,@("V13", $null, $null, 7408, 7410, $null)
$data += ,@("V14", $null, $null, 7408, 7410, $null)
$data += ,@("V15", $null, $null, 7400, 7400, $null)
$data += ,@("V16", $null, $null, 7432, 7411, $null)
$data += ,@("V17", $null, $null, 7432, 7430, $null)
$data += ,@("V18", $null, $null, $null, 7430, $null)
$data += ,@("V19", $null, $null, $null, 7432, $null)
$data += ,@("V20", $null, $null, $null, 7400, $null)
$data += ,@("V21", $null, $null, $null, 7430, $null)
$data += ,@("V22", $null, $null, $null, 7432, $null)
$data += ,@("V23", $null, $null, $null, 7411, $null)
$data += ,@("V24", $null, $null, $null, 7430, $null)
$data += ,@("V25", $null, $null, $null, 7410, $null)
$data += ,@("V26", $null, $null, $null, 7410, $null)

$cols = @("B", "C", "D", "E", "F")
$rowNum = 27
foreach ($row in $data) {
    $ws.Range("A$rowNum").Value = $row[0]
    for ($i = 0; $i -lt 5; $i++) {
        $val = $row[$i + 1]
        if ($null -ne $val) {
            $colLetter = $cols[$i]
            $ws.Range("$colLetter$rowNum").Value = $val
        }
    }
    $rowNum++
}

# Restore the sheet as the sole, active tab with the on-screen selection
# last left at G32 (scrolled so row 17 is at the top).
$ws.Activate()
[void]$ws.Range("G32").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
